# Slide 15 ("總計") holds a summary table of 有效樣本數 (valid sample
# counts) per subject for 心跳 (heart rate) and 呼吸 (breathing rate).
# Corrected counts for the "hung-wei" row and the derived total/accuracy
# figures in the 呼吸 (breath) column (plus the 心跳 total).
#
# Table columns: 1 = subject name, 2 = 心跳 (heart), 3 = 呼吸 (breath)
#   Row 3  (hung-wei)   : col 3 10        -> 9
#   Row 15 (total)      : col 2 114       -> 113
#                          col 3 63        -> 59
#   Row 16 (Accuracy/%) : col 3 37/58.7%  -> 36/61%
#
# Note: Table.Rows.Item(r).Cells.Item(c).Shape only resolves to the first
# cell of row r in this host, so cells must be addressed through
# Table.Cell(row, column), which correctly resolves per-column shapes.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)

$tbl = $null
foreach ($shp in $s.Shapes) {
    if ($shp.HasTable) {
        $tbl = $shp.Table
        break
    }
}

$tbl.Cell(3, 3).Shape.TextFrame.TextRange.Text = "9"
$tbl.Cell(15, 2).Shape.TextFrame.TextRange.Text = "113"
$tbl.Cell(15, 3).Shape.TextFrame.TextRange.Text = "59"
$tbl.Cell(16, 3).Shape.TextFrame.TextRange.Text = "36/61%"
